# Generate Report for handoff
# Renames the source markdown file that was handed off (ff79a594... -> fb427a9c...),
# records that its handoff transform failed, regenerates fresh target .xlf hand-off
# files (new uuid/hash/timestamps), and keeps the previously-ignored
# ".localization-config" file as the last row on every sheet.

$wb = $excel.ActiveWorkbook

$oldMdName = "ff79a594-5c11-46b5-91b9-8eddece4be7c.md"
$newMdName = "fb427a9c-67d1-49c4-93d8-81f46a664ad9.md"
$newMdName2 = "5e7ad45d-d58a-4a79-826a-ee964522209f.md"
$failedStatus = "Handoff transform failed"

$oldHashZh = "ff79a594-5c11-46b5-91b9-8eddece4be7c.43d04fe0cd66395c62bac4d6b6c1d9e35997abb5.zh-cn.xlf"
$newHashZh = "fb427a9c-67d1-49c4-93d8-81f46a664ad9.e1796fb6e7146e9356c516e009cdd6521abd467a.zh-cn.xlf"
$newHandoffZh = "2016-01-11 17:06:34"

$oldHashDe = "ff79a594-5c11-46b5-91b9-8eddece4be7c.43d04fe0cd66395c62bac4d6b6c1d9e35997abb5.de-de.xlf"
$newHashDe = "fb427a9c-67d1-49c4-93d8-81f46a664ad9.e1796fb6e7146e9356c516e009cdd6521abd467a.de-de.xlf"
$newHandoffDe = "2016-01-11 17:06:56"

$epoch = "0001-01-01 00:00:00"
$ignored = "Ignored"

function Rewrite-GithubUrl($oldUrl, $oldToken, $newToken) {
    return $oldUrl.Replace($oldToken, $newToken)
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# remember the original hyperlink targets before we touch anything
$ws1Link2Addr = $null
$ws1Link3Addr = $null
foreach ($h in $ws1.Hyperlinks) {
    if ($h.Range.Address() -eq "`$A`$2") { $ws1Link2Addr = $h.Address }
    if ($h.Range.Address() -eq "`$A`$3") { $ws1Link3Addr = $h.Address }
}

# push the ".localization-config" row down from row 3 to row 4, and free up
# row 3 for the newly generated "transform failed" entry
$ws1.Rows(3).Copy()
$ws1.Rows(3).Insert()

# row 2 : renamed source file, status unchanged
$ws1.Range("A2").Value = $newMdName

# row 3 : brand-new entry for the file whose handoff transform failed
$ws1.Range("A3").Value = $newMdName2
$ws1.Range("B3").Value = $failedStatus
$ws1.Range("C3").Value = $failedStatus

# row 4 already contains ".localization-config" / "Not to be localized"
# (duplicated by the row insert above), nothing else to change there.

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), (Rewrite-GithubUrl $ws1Link2Addr "ff79a594-5c11-46b5-91b9-8eddece4be7c" "fb427a9c-67d1-49c4-93d8-81f46a664ad9"), "", "", $newMdName)
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($ws1Link2Addr.Replace("ff79a594-5c11-46b5-91b9-8eddece4be7c", "5e7ad45d-d58a-4a79-826a-ee964522209f")), "", "", $newMdName2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $ws1Link3Addr, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2Link2Addr = $null
$ws2Link3Addr = $null
$ws2Link4Addr = $null
foreach ($h in $ws2.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") { $ws2Link2Addr = $h.Address }
    if ($addr -eq "`$C`$2") { $ws2Link3Addr = $h.Address }
    if ($addr -eq "`$A`$3") { $ws2Link4Addr = $h.Address }
}

$ws2.Rows(3).Copy()
$ws2.Rows(3).Insert()
$ws2.Range("C3").Clear()

# row 2 : renamed source file + refreshed hand-off artifacts
$ws2.Range("A2").Value = $newMdName
$ws2.Range("C2").Value = $newHashZh
$ws2.Range("D2").Value = $newHandoffZh

# row 3 : brand-new entry for the file whose handoff transform failed
$ws2.Range("A3").Value = $newMdName2
$ws2.Range("B3").Value = $failedStatus
$ws2.Range("D3").Value = $epoch
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = $ignored

# row 4 keeps the ".localization-config" / "Not to be localized" values that
# were duplicated by the row insert above.

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), (Rewrite-GithubUrl $ws2Link2Addr "ff79a594-5c11-46b5-91b9-8eddece4be7c" "fb427a9c-67d1-49c4-93d8-81f46a664ad9"), "", "", $newMdName)
$ws2.Hyperlinks.Add($ws2.Range("C2"), (Rewrite-GithubUrl (Rewrite-GithubUrl $ws2Link3Addr "ff79a594-5c11-46b5-91b9-8eddece4be7c" "fb427a9c-67d1-49c4-93d8-81f46a664ad9") "43d04fe0cd66395c62bac4d6b6c1d9e35997abb5" "e1796fb6e7146e9356c516e009cdd6521abd467a"), "", "", $newHashZh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($ws2Link2Addr.Replace("ff79a594-5c11-46b5-91b9-8eddece4be7c", "5e7ad45d-d58a-4a79-826a-ee964522209f")), "", "", $newMdName2)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $ws2Link4Addr, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3Link2Addr = $null
$ws3Link3Addr = $null
$ws3Link4Addr = $null
foreach ($h in $ws3.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") { $ws3Link2Addr = $h.Address }
    if ($addr -eq "`$C`$2") { $ws3Link3Addr = $h.Address }
    if ($addr -eq "`$A`$3") { $ws3Link4Addr = $h.Address }
}

$ws3.Rows(3).Copy()
$ws3.Rows(3).Insert()
$ws3.Range("C3").Clear()

# row 2 : renamed source file + refreshed hand-off artifacts
$ws3.Range("A2").Value = $newMdName
$ws3.Range("C2").Value = $newHashDe
$ws3.Range("D2").Value = $newHandoffDe

# row 3 : brand-new entry for the file whose handoff transform failed
$ws3.Range("A3").Value = $newMdName2
$ws3.Range("B3").Value = $failedStatus
$ws3.Range("D3").Value = $epoch
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = $ignored

# row 4 keeps the ".localization-config" / "Not to be localized" values that
# were duplicated by the row insert above.

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), (Rewrite-GithubUrl $ws3Link2Addr "ff79a594-5c11-46b5-91b9-8eddece4be7c" "fb427a9c-67d1-49c4-93d8-81f46a664ad9"), "", "", $newMdName)
$ws3.Hyperlinks.Add($ws3.Range("C2"), (Rewrite-GithubUrl (Rewrite-GithubUrl $ws3Link3Addr "ff79a594-5c11-46b5-91b9-8eddece4be7c" "fb427a9c-67d1-49c4-93d8-81f46a664ad9") "43d04fe0cd66395c62bac4d6b6c1d9e35997abb5" "e1796fb6e7146e9356c516e009cdd6521abd467a"), "", "", $newHashDe)
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($ws3Link2Addr.Replace("ff79a594-5c11-46b5-91b9-8eddece4be7c", "5e7ad45d-d58a-4a79-826a-ee964522209f")), "", "", $newMdName2)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $ws3Link4Addr, "", "", ".localization-config")

Write-Host "Report regenerated for handoff"
